$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (2-9); keep header row (and its shared strings) intact
$ws.Range("A2:T9").ClearContents()

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Angpt1"
$ws.Cells.Item(2, 3).Value = "Itgb1"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.09342833333333334
$ws.Cells.Item(2, 8).Value = 0.280285
$ws.Cells.Item(2, 9).Value = 0.004971127078137211
$ws.Cells.Item(2, 10).Value = 0.004971127078137211
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 98.946724
$ws.Cells.Item(2, 14).Value = 296.840172
$ws.Cells.Item(2, 15).Value = 0.2098009692989996
$ws.Cells.Item(2, 16).Value = 0.2098009692989996
$ws.Cells.Item(2, 17).Value = 9.244427512113333
$ws.Cells.Item(2, 18).Value = 83.19984760902
$ws.Cells.Item(2, 19).Value = 0.00104294727950169
$ws.Cells.Item(2, 20).Value = 0.00104294727950169

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Angpt1"
$ws.Cells.Item(3, 3).Value = "Itgb1"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.09342833333333334
$ws.Cells.Item(3, 8).Value = 0.280285
$ws.Cells.Item(3, 9).Value = 0.004971127078137211
$ws.Cells.Item(3, 10).Value = 0.004971127078137211
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 163.0062356666667
$ws.Cells.Item(3, 14).Value = 489.018707
$ws.Cells.Item(3, 15).Value = 0.345629090707923
$ws.Cells.Item(3, 16).Value = 0.3456290907079231
$ws.Cells.Item(3, 17).Value = 15.22940092127722
$ws.Cells.Item(3, 18).Value = 137.064608291495
$ws.Cells.Item(3, 19).Value = 0.001718166131810098
$ws.Cells.Item(3, 20).Value = 0.001718166131810099

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Angpt1"
$ws.Cells.Item(4, 3).Value = "Itgb1"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.09342833333333334
$ws.Cells.Item(4, 8).Value = 0.280285
$ws.Cells.Item(4, 9).Value = 0.004971127078137211
$ws.Cells.Item(4, 10).Value = 0.004971127078137211
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 65.39610666666668
$ws.Cells.Item(4, 14).Value = 196.18832
$ws.Cells.Item(4, 15).Value = 0.1386621609326595
$ws.Cells.Item(4, 16).Value = 0.1386621609326595
$ws.Cells.Item(4, 17).Value = 6.109849252355557
$ws.Cells.Item(4, 18).Value = 54.98864327120001
$ws.Cells.Item(4, 19).Value = 0.0006893072229253631
$ws.Cells.Item(4, 20).Value = 0.0006893072229253633

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Angpt1"
$ws.Cells.Item(5, 3).Value = "Itgb1"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.09342833333333334
$ws.Cells.Item(5, 8).Value = 0.280285
$ws.Cells.Item(5, 9).Value = 0.004971127078137211
$ws.Cells.Item(5, 10).Value = 0.004971127078137211
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 144.2727966666667
$ws.Cells.Item(5, 14).Value = 432.81839
$ws.Cells.Item(5, 15).Value = 0.3059077790604178
$ws.Cells.Item(5, 16).Value = 0.3059077790604179
$ws.Cells.Item(5, 17).Value = 13.47916693790556
$ws.Cells.Item(5, 18).Value = 121.31250244115
$ws.Cells.Item(5, 19).Value = 0.001520706443900058
$ws.Cells.Item(5, 20).Value = 0.001520706443900059

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Angpt1"
$ws.Cells.Item(6, 3).Value = "Itgb1"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 14.289121
$ws.Cells.Item(6, 8).Value = 42.867363
$ws.Cells.Item(6, 9).Value = 0.760294375288143
$ws.Cells.Item(6, 10).Value = 0.760294375288143
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 98.946724
$ws.Cells.Item(6, 14).Value = 296.840172
$ws.Cells.Item(6, 15).Value = 0.2098009692989996
$ws.Cells.Item(6, 16).Value = 0.2098009692989996
$ws.Cells.Item(6, 17).Value = 1413.861711789604
$ws.Cells.Item(6, 18).Value = 12724.75540610643
$ws.Cells.Item(6, 19).Value = 0.1595104968880298
$ws.Cells.Item(6, 20).Value = 0.1595104968880298

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Angpt1"
$ws.Cells.Item(7, 3).Value = "Itgb1"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 14.289121
$ws.Cells.Item(7, 8).Value = 42.867363
$ws.Cells.Item(7, 9).Value = 0.760294375288143
$ws.Cells.Item(7, 10).Value = 0.760294375288143
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 163.0062356666667
$ws.Cells.Item(7, 14).Value = 489.018707
$ws.Cells.Item(7, 15).Value = 0.345629090707923
$ws.Cells.Item(7, 16).Value = 0.3456290907079231
$ws.Cells.Item(7, 17).Value = 2329.215825195516
$ws.Cells.Item(7, 18).Value = 20962.94242675964
$ws.Cells.Item(7, 19).Value = 0.2627798536011893
$ws.Cells.Item(7, 20).Value = 0.2627798536011893

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Angpt1"
$ws.Cells.Item(8, 3).Value = "Itgb1"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 14.289121
$ws.Cells.Item(8, 8).Value = 42.867363
$ws.Cells.Item(8, 9).Value = 0.760294375288143
$ws.Cells.Item(8, 10).Value = 0.760294375288143
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 65.39610666666668
$ws.Cells.Item(8, 14).Value = 196.18832
$ws.Cells.Item(8, 15).Value = 0.1386621609326595
$ws.Cells.Item(8, 16).Value = 0.1386621609326595
$ws.Cells.Item(8, 17).Value = 934.4528810889069
$ws.Cells.Item(8, 18).Value = 8410.075929800161
$ws.Cells.Item(8, 19).Value = 0.1054240610224003
$ws.Cells.Item(8, 20).Value = 0.1054240610224003

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Angpt1"
$ws.Cells.Item(9, 3).Value = "Itgb1"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 14.289121
$ws.Cells.Item(9, 8).Value = 42.867363
$ws.Cells.Item(9, 9).Value = 0.760294375288143
$ws.Cells.Item(9, 10).Value = 0.760294375288143
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 144.2727966666667
$ws.Cells.Item(9, 14).Value = 432.81839
$ws.Cells.Item(9, 15).Value = 0.3059077790604178
$ws.Cells.Item(9, 16).Value = 0.3059077790604179
$ws.Cells.Item(9, 17).Value = 2061.531448578397
$ws.Cells.Item(9, 18).Value = 18553.78303720557
$ws.Cells.Item(9, 19).Value = 0.2325799637765237
$ws.Cells.Item(9, 20).Value = 0.2325799637765237

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Angpt1"
$ws.Cells.Item(10, 3).Value = "Itgb1"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 4.411646
$ws.Cells.Item(10, 8).Value = 13.234938
$ws.Cells.Item(10, 9).Value = 0.2347344976337197
$ws.Cells.Item(10, 10).Value = 0.2347344976337198
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 98.946724
$ws.Cells.Item(10, 14).Value = 296.840172
$ws.Cells.Item(10, 15).Value = 0.2098009692989996
$ws.Cells.Item(10, 16).Value = 0.2098009692989996
$ws.Cells.Item(10, 17).Value = 436.517919147704
$ws.Cells.Item(10, 18).Value = 3928.661272329336
$ws.Cells.Item(10, 19).Value = 0.04924752513146812
$ws.Cells.Item(10, 20).Value = 0.04924752513146814

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Angpt1"
$ws.Cells.Item(11, 3).Value = "Itgb1"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 4.411646
$ws.Cells.Item(11, 8).Value = 13.234938
$ws.Cells.Item(11, 9).Value = 0.2347344976337197
$ws.Cells.Item(11, 10).Value = 0.2347344976337198
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 163.0062356666667
$ws.Cells.Item(11, 14).Value = 489.018707
$ws.Cells.Item(11, 15).Value = 0.345629090707923
$ws.Cells.Item(11, 16).Value = 0.3456290907079231
$ws.Cells.Item(11, 17).Value = 719.1258075539074
$ws.Cells.Item(11, 18).Value = 6472.132267985166
$ws.Cells.Item(11, 19).Value = 0.08113107097492366
$ws.Cells.Item(11, 20).Value = 0.0811310709749237

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Angpt1"
$ws.Cells.Item(12, 3).Value = "Itgb1"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 4.411646
$ws.Cells.Item(12, 8).Value = 13.234938
$ws.Cells.Item(12, 9).Value = 0.2347344976337197
$ws.Cells.Item(12, 10).Value = 0.2347344976337198
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 65.39610666666668
$ws.Cells.Item(12, 14).Value = 196.18832
$ws.Cells.Item(12, 15).Value = 0.1386621609326595
$ws.Cells.Item(12, 16).Value = 0.1386621609326595
$ws.Cells.Item(12, 17).Value = 288.5044723915734
$ws.Cells.Item(12, 18).Value = 2596.540251524161
$ws.Cells.Item(12, 19).Value = 0.03254879268733381
$ws.Cells.Item(12, 20).Value = 0.03254879268733383

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Angpt1"
$ws.Cells.Item(13, 3).Value = "Itgb1"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 4.411646
$ws.Cells.Item(13, 8).Value = 13.234938
$ws.Cells.Item(13, 9).Value = 0.2347344976337197
$ws.Cells.Item(13, 10).Value = 0.2347344976337198
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 144.2727966666667
$ws.Cells.Item(13, 14).Value = 432.81839
$ws.Cells.Item(13, 15).Value = 0.3059077790604178
$ws.Cells.Item(13, 16).Value = 0.3059077790604179
$ws.Cells.Item(13, 17).Value = 636.4805063233133
$ws.Cells.Item(13, 18).Value = 5728.324556909821
$ws.Cells.Item(13, 19).Value = 0.07180710883999411
$ws.Cells.Item(13, 20).Value = 0.07180710883999414

Write-Output "done"
